$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 92, shifting existing rows 92:103 down to 93:104.
$ws.Rows("92:92").Insert()

# Populate the newly inserted row 92 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T repeat the same market/product info
# as the rest of this block; D,M,N,O,P,S carry the new observation.
$ws.Cells.Item(92, 1).Value = 10
$ws.Cells.Item(92, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(92, 3).Value = "La Araucanía"
$ws.Cells.Item(92, 4).Value = 45124
$ws.Cells.Item(92, 5).Value = 9
$ws.Cells.Item(92, 6).Value = "Fruta"
$ws.Cells.Item(92, 7).Value = 100108
$ws.Cells.Item(92, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(92, 9).Value = 100108007
$ws.Cells.Item(92, 10).Value = "Coco"
$ws.Cells.Item(92, 11).Value = "Sin especificar"
$ws.Cells.Item(92, 12).Value = "Primera"
$ws.Cells.Item(92, 13).Value = 15
$ws.Cells.Item(92, 14).Value = 36000
$ws.Cells.Item(92, 15).Value = 36000
$ws.Cells.Item(92, 16).Value = 36000
$ws.Cells.Item(92, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(92, 18).Value = "Perú"
$ws.Cells.Item(92, 19).Value = 1800
$ws.Cells.Item(92, 20).Value = 20
